$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E11: "Pendiente ADM" -> " 01229548" (leading space, keep as text) ---
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = " 01229548"

# --- New row 18 ---
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "4757 "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "12/11/2025"
$ws.Range("C18").Value = "GARAY, JUAN DE AV. 819"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Pendiente ADM"
$ws.Range("F18").Value = "Optical Power"
$ws.Range("G18").Value = "Pendiente"
$ws.Range("H18").Value = "tendido bajo"
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = '{"direccionesNormalizadas": [{"altura": 819, "cod_calle": 7026, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.376986", "y": "-34.625210"}, "direccion": "GARAY, JUAN DE AV. 819, CABA", "nombre_calle": "GARAY, JUAN DE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K18").Value = -58.376986
$ws.Range("L18").Value = -34.62521
$ws.Range("M18").Value = "San Telmo"
$ws.Range("N18").Value = "Capital Sur"

# --- New row 19 ---
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "4756 "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "12/11/2025"
$ws.Range("C19").Value = "GARAY, JUAN DE AV. 799"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "Pendiente ADM"
$ws.Range("F19").Value = "Optical Power"
$ws.Range("G19").Value = "Pendiente"
$ws.Range("H19").Value = "tendido bajo"
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = '{"direccionesNormalizadas": [{"altura": 799, "cod_calle": 7026, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.376455", "y": "-34.624886"}, "direccion": "GARAY, JUAN DE AV. 799, CABA", "nombre_calle": "GARAY, JUAN DE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K19").Value = -58.376455
$ws.Range("L19").Value = -34.624886
$ws.Range("M19").Value = "San Telmo"
$ws.Range("N19").Value = "Capital Sur"
